$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.986.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.148.54'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.07%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.141.37'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.536'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.72%  '
$ws.Range('E10').Value = '  +17.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.74'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.471'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('E13').Value = '  +7.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.53%  '
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.670.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.938.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.144.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '471.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.735'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.32%  '
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.30%  '
$ws.Range('E28').Value = '  +4.76%  '
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  +3.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0875'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.39%  '
$ws.Range('E36').Value = '  +3.61%  '
$ws.Range('E37').Value = '  +15.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.95'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '450.09'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.39%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  +5.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.926.90'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.284'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.47%  '
$ws.Range('E45').Value = '  +5.03%  '
$ws.Range('E46').Value = '  +6.19%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '124.63'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.112'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.84%  '
